$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General->Text) on column D cells being updated so that numeric-looking
# values (e.g. "415.06", "1.00") are stored as text, matching the original inlineStr type
# instead of being auto-converted to numbers by Excel.
$dCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D49", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values
$ws.Range('D2').Value = '62.419.15'
$ws.Range('E2').Value = '  +9.34%  '
$ws.Range('D3').Value = '3.368.81'
$ws.Range('E3').Value = '  +3.72%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '415.06'
$ws.Range('E5').Value = '  +4.73%  '
$ws.Range('D6').Value = '116.54'
$ws.Range('E6').Value = '  +7.58%  '
$ws.Range('D7').Value = '3.359.19'
$ws.Range('E7').Value = '  +3.56%  '
$ws.Range('D8').Value = '0.576'
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '0.631'
$ws.Range('E10').Value = '  +0.93%  '
$ws.Range('D11').Value = '0.118'
$ws.Range('E11').Value = '  +19.78%  '
$ws.Range('D12').Value = '40.09'
$ws.Range('E12').Value = '  +2.10%  '
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').Value = '3.885.59'
$ws.Range('E14').Value = '  +3.41%  '
$ws.Range('D15').Value = '8.34'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').Value = '19.31'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '3.356.48'
$ws.Range('E17').Value = '  +3.16%  '
$ws.Range('D18').Value = '62.000.76'
$ws.Range('E18').Value = '  +8.98%  '
$ws.Range('E19').Value = '  -2.06%  '
$ws.Range('D20').Value = '10.87'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').Value = '0.0000118'
$ws.Range('E21').Value = '  +8.91%  '
$ws.Range('D22').Value = '3.34'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '12.56'
$ws.Range('E23').Value = '  -4.07%  '
$ws.Range('D24').Value = '295.99'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').Value = '74.82'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').Value = '3.14'
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('D27').Value = '29.45'
$ws.Range('E27').Value = '  +4.75%  '
$ws.Range('D28').Value = '7.95'
$ws.Range('E28').Value = '  +9.90%  '
$ws.Range('E29').Value = '  +3.70%  '
$ws.Range('D30').Value = '4.27'
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('D31').Value = '7.63'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('D32').Value = '43.17'
$ws.Range('E32').Value = '  +8.16%  '
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').Value = '2.56'
$ws.Range('E34').Value = '  +20.29%  '
$ws.Range('B35').Value = 'Cosmos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D35').Value = '11.47'
$ws.Range('E35').Value = '  +2.08%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('D38').Value = '52.31'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').Value = '3.14'
$ws.Range('E39').Value = '  +6.73%  '
$ws.Range('D40').Value = '0.995'
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('D41').Value = '3.46'
$ws.Range('E41').Value = '  -0.82%  '
$ws.Range('D42').Value = '133.42'
$ws.Range('E42').Value = '  -3.25%  '
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('D44').Value = '1.91'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').Value = '0.287'
$ws.Range('E45').Value = '  +2.45%  '
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('E47').Value = '  -3.43%  '
$ws.Range('E48').Value = '  -3.46%  '
$ws.Range('D49').Value = '2.171.64'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('D50').Value = '21.22'
$ws.Range('E50').Value = '  -4.46%  '
$ws.Range('D51').Value = '3.688.31'
$ws.Range('E51').Value = '  +3.29%  '
